# Add my API key
# - "Simple Fields" sheet: D2 telephone-number becomes the API key text
# - "Simple Fields - Formatted" sheet: D2 becomes a CSV-ish multi-line blob,
#   and needs word-wrap turned on so the embedded newline renders.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Simple Fields")
$ws1.Range("D2").Value = "01152019for"

$ws2 = $wb.Worksheets.Item("Simple Fields - Formatted")
$ws2.Range("D2").Value = "Key,Value`n`"Value`",`"`""
$ws2.Range("D2").WrapText = $true
